$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Replace "Souhvězdí Blíženci." with "Souhvězdí Blíženců." everywhere in the document
# (wdReplaceAll = 2)
$find.Execute("zobrazujíSouhvězdí Blíženci.", $true, $false, $false, $false, $false, `
              $true, 1, $false, "zobrazujíSouhvězdí Blíženců.", 2) | Out-Null
